# Apply updated "want-to-go" counts / min-price figures to both the
# "展览" sheet and the combined "全部类型" sheet (which mirrors the same
# rows, shifted down by one).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8
$ws1.Range("F3").Value = 12898
$ws1.Range("G3").Value = 60
$ws1.Range("F5").Value = 88
$ws1.Range("F6").Value = 72
$ws1.Range("F10").Value = 12826
$ws1.Range("F11").Value = 280
$ws1.Range("F12").Value = 36
$ws1.Range("F13").Value = 8683
$ws1.Range("F14").Value = 7678
$ws1.Range("F16").Value = 108
$ws1.Range("F18").Value = 125
$ws1.Range("F24").Value = 17

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 8
$ws4.Range("F4").Value = 12898
$ws4.Range("G4").Value = 60
$ws4.Range("F6").Value = 88
$ws4.Range("F7").Value = 72
$ws4.Range("F11").Value = 12826
$ws4.Range("F12").Value = 280
$ws4.Range("F13").Value = 36
$ws4.Range("F14").Value = 8683
$ws4.Range("F15").Value = 7678
$ws4.Range("F17").Value = 108
$ws4.Range("F19").Value = 125
$ws4.Range("F26").Value = 17
